$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 12424
$ws.Range("F3").Value = 6963
$ws.Range("F11").Value = 125
$ws.Range("F12").Value = 327
$ws.Range("F13").Value = 981
$ws.Range("F14").Value = 3709
$ws.Range("F15").Value = 65
$ws.Range("F18").Value = 218
$ws.Range("F19").Value = 349
$ws.Range("F22").Value = 290
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 95
$ws.Range("F25").Value = 342
$ws.Range("F26").Value = 5142
$ws.Range("F28").Value = 1367
$ws.Range("F29").Value = 277
$ws.Range("F30").Value = 903
$ws.Range("F31").Value = 1296

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 3733

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9196
$ws.Range("F3").Value = 542
$ws.Range("F4").Value = 1934

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9196
$ws.Range("F3").Value = 542
$ws.Range("F4").Value = 1934
$ws.Range("F5").Value = 12424
$ws.Range("F6").Value = 6963
$ws.Range("F8").Value = 3733
$ws.Range("F16").Value = 125
$ws.Range("F17").Value = 327
$ws.Range("F18").Value = 981
$ws.Range("F19").Value = 3709
$ws.Range("F20").Value = 65
$ws.Range("F22").Value = 218
$ws.Range("F23").Value = 349
$ws.Range("F26").Value = 290
$ws.Range("F27").Value = 29
$ws.Range("F32").Value = 342
$ws.Range("F33").Value = 5142
$ws.Range("F35").Value = 1367
$ws.Range("F38").Value = 277
$ws.Range("F40").Value = 904
$ws.Range("F41").Value = 1296
